# ----------------------------------------------------------------------------
# Rename the worksheet ("Sheet" -> "Sheet1")
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# ----------------------------------------------------------------------------
# Header row (row 1): values first, formatting afterwards
# ----------------------------------------------------------------------------
$ws.Range("A1").Value = 'Date'
$ws.Range("B1").Value = 'Model Name'
$ws.Range("C1").Value = 'Exact Precision (Micro Avg)'
$ws.Range("D1").Value = 'Exact Recall (Micro Avg)'
$ws.Range("E1").Value = 'Exact F1 Score (Micro Avg)'
$ws.Range("F1").Value = 'Exact Precision (Macro Avg)'
$ws.Range("G1").Value = 'Exact Recall (Macro Avg)'
$ws.Range("H1").Value = 'Exact F1 Score (Macro Avg)'
$ws.Range("I1").Value = 'Exact Precision (Weighted Avg)'
$ws.Range("J1").Value = 'Exact Recall (Weighted Avg)'
$ws.Range("K1").Value = 'Exact F1 Score (Weighted Avg)'
$ws.Range("L1").Value = 'Partial Precision'
$ws.Range("M1").Value = 'Partial Recall'
$ws.Range("N1").Value = 'Partial F1 Score'
$ws.Range("O1").Value = 'Partial TP'
$ws.Range("P1").Value = 'Partial FP'
$ws.Range("Q1").Value = 'Partial FN'
$ws.Range("R1").Value = 'Support'
$ws.Range("S1").Value = 'Accuracy'
$ws.Range("T1").Value = 'Result Link'
$ws.Range("U1").Value = 'Stats Link'
$ws.Range("V1").Value = 'No of GPU Used'
$ws.Range("W1").Value = 'Power Consumption'
$ws.Range("X1").Value = 'Unnamed: 23'

# --- Data row 2 ---
$ws.Range("B2").Value = 'Qwen2.5-32B-Instruct'
$ws.Range("C2").Value = 0.5454545454545454
$ws.Range("D2").Value = 0.3434343434343434
$ws.Range("E2").Value = 0.4214876033057851
$ws.Range("F2").Value = 0.2576827202404189
$ws.Range("G2").Value = 0.1677352029058214
$ws.Range("H2").Value = 0.1993802158805046
$ws.Range("I2").Value = 0.5680544317482078
$ws.Range("J2").Value = 0.3434343434343434
$ws.Range("K2").Value = 0.4213704874378936
$ws.Range("L2").Value = 0.654054054054054
$ws.Range("M2").Value = 0.4087837837837838
$ws.Range("N2").Value = 0.5031185031185031
$ws.Range("O2").Value = 121
$ws.Range("P2").Value = 64
$ws.Range("Q2").Value = 175
$ws.Range("R2").Value = 297
$ws.Range("S2").Value = 0.9605751947273816
$ws.Range("T2").Value = '/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_5_shot.txt'
$ws.Range("U2").Value = '/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_5_shot.txt'
$ws.Range("V2").Value = '4 MLGPU'
$ws.Range("W2").Value = '0.060 kWh'
$ws.Range("X2").Value = 2141

# --- Data row 3 ---
$ws.Range("B3").Value = 'Qwen2.5-32B-Instruct'
$ws.Range("C3").Value = 0.4937759336099585
$ws.Range("D3").Value = 0.4006734006734007
$ws.Range("E3").Value = 0.4423791821561339
$ws.Range("F3").Value = 0.6011950196038596
$ws.Range("G3").Value = 0.4036792880542693
$ws.Range("H3").Value = 0.4598232623324498
$ws.Range("I3").Value = 0.5919217243387112
$ws.Range("J3").Value = 0.4006734006734007
$ws.Range("K3").Value = 0.457523533725503
$ws.Range("L3").Value = 0.5560165975103735
$ws.Range("M3").Value = 0.4527027027027027
$ws.Range("N3").Value = 0.4990689013035382
$ws.Range("O3").Value = 134
$ws.Range("P3").Value = 107
$ws.Range("Q3").Value = 162
$ws.Range("R3").Value = 297
$ws.Range("S3").Value = 0.9529059316956261
$ws.Range("T3").Value = '/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_5_shot.txt'
$ws.Range("U3").Value = '/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_5_shot.txt'
$ws.Range("V3").Value = '4 MLGPU'
$ws.Range("W3").Value = '0.079 kWh'

# ----------------------------------------------------------------------------
# Force the "Date" column values (A2, A3) to remain plain text rather than
# being auto-converted to Excel date serials. We build the text format on a
# scratch cell far outside the used range, copy its format onto A2:A3, then
# assign the literal strings, and finally wipe the scratch cell completely
# (value + format) so it leaves no trace in the saved workbook.
# ----------------------------------------------------------------------------
$dateFmtTemplate = $ws.Range("AA1")
$dateFmtTemplate.NumberFormat = "@"
$dateFmtTemplate.Value = "x"
$dateFmtTemplate.Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$dateFmtTemplate.Clear()

$ws.Range("A2").Value = "09/11/2025"
$ws.Range("A3").Value = "09/12/2025"

# ----------------------------------------------------------------------------
# X3 must exist as an empty (but present) cell so the sheet dimension extends
# to column X on row 3. Touching a border property with the "no line" value
# materializes the cell without allocating a new distinct style.
# ----------------------------------------------------------------------------
$ws.Range("X3").Borders.LineStyle = 0

# ----------------------------------------------------------------------------
# Header formatting: bold font, thin box border, centered horizontally and
# top-aligned vertically. Built once on a scratch cell and then copied onto
# the header range in a single PasteSpecial so only one extra style entry is
# produced (matching the two-entry cellXfs table of the target workbook).
# ----------------------------------------------------------------------------
$headerFmtTemplate = $ws.Range("AA2")
$headerFmtTemplate.Value = "x"
$headerFmtTemplate.Font.Bold = $true
$headerFmtTemplate.Borders.LineStyle = 1
$headerFmtTemplate.HorizontalAlignment = -4108
$headerFmtTemplate.VerticalAlignment = -4160

$headerRange = $ws.Range("A1:X1")
$headerFmtTemplate.Copy()
$headerRange.PasteSpecial(-4122)
$headerFmtTemplate.Clear()

Write-Host "edit complete"
